{"js": "// Replace the date line and each \"AAA\u00f7B=\" division expression in the\n// worksheet with the updated values from the commit. Every old value is\n// unique within the document, so a simple search+replace per pair is safe.\nconst replacements = [\n  [\"2025-07-17 Thursday\", \"2025-07-18 Friday\"],\n  [\"396\u00f74=\", \"191\u00f74=\"],\n  [\"944\u00f73=\", \"749\u00f78=\"],\n  [\"930\u00f79=\", \"226\u00f74=\"],\n  [\"147\u00f73=\", \"379\u00f73=\"],\n  [\"962\u00f76=\", \"118\u00f72=\"],\n  [\"671\u00f78=\", \"374\u00f75=\"],\n  [\"906\u00f73=\", \"797\u00f77=\"],\n  [\"555\u00f75=\", \"727\u00f76=\"],\n  [\"316\u00f72=\", \"385\u00f77=\"],\n  [\"556\u00f76=\", \"444\u00f78=\"],\n  [\"185\u00f76=\", \"847\u00f79=\"],\n  [\"504\u00f78=\", \"743\u00f73=\"],\n  [\"375\u00f72=\", \"718\u00f72=\"],\n  [\"422\u00f75=\", \"581\u00f79=\"],\n  [\"526\u00f77=\", \"345\u00f76=\"],\n  [\"254\u00f72=\", \"446\u00f79=\"],\n  [\"490\u00f74=\", \"985\u00f75=\"],\n  [\"514\u00f72=\", \"400\u00f76=\"],\n  [\"199\u00f72=\", \"427\u00f75=\"],\n  [\"200\u00f72=\", \"987\u00f77=\"],\n  [\"531\u00f75=\", \"519\u00f75=\"],\n  [\"152\u00f77=\", \"223\u00f73=\"],\n  [\"670\u00f78=\", \"755\u00f73=\"],\n  [\"545\u00f74=\", \"439\u00f75=\"],\n  [\"723\u00f72=\", \"688\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"AAA\u00f7B=\" division expression with the\n# updated values from the commit. Every \"Old\" value is unique within the\n# document, so Find/Replace All (wrap = none needed, single hit each) is\n# safe to run independently for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"2025-07-17 Thursday\"; New = \"2025-07-18 Friday\" },\n  @{ Old = \"396\u00f74=\"; New = \"191\u00f74=\" },\n  @{ Old = \"944\u00f73=\"; New = \"749\u00f78=\" },\n  @{ Old = \"930\u00f79=\"; New = \"226\u00f74=\" },\n  @{ Old = \"147\u00f73=\"; New = \"379\u00f73=\" },\n  @{ Old = \"962\u00f76=\"; New = \"118\u00f72=\" },\n  @{ Old = \"671\u00f78=\"; New = \"374\u00f75=\" },\n  @{ Old = \"906\u00f73=\"; New = \"797\u00f77=\" },\n  @{ Old = \"555\u00f75=\"; New = \"727\u00f76=\" },\n  @{ Old = \"316\u00f72=\"; New = \"385\u00f77=\" },\n  @{ Old = \"556\u00f76=\"; New = \"444\u00f78=\" },\n  @{ Old = \"185\u00f76=\"; New = \"847\u00f79=\" },\n  @{ Old = \"504\u00f78=\"; New = \"743\u00f73=\" },\n  @{ Old = \"375\u00f72=\"; New = \"718\u00f72=\" },\n  @{ Old = \"422\u00f75=\"; New = \"581\u00f79=\" },\n  @{ Old = \"526\u00f77=\"; New = \"345\u00f76=\" },\n  @{ Old = \"254\u00f72=\"; New = \"446\u00f79=\" },\n  @{ Old = \"490\u00f74=\"; New = \"985\u00f75=\" },\n  @{ Old = \"514\u00f72=\"; New = \"400\u00f76=\" },\n  @{ Old = \"199\u00f72=\"; New = \"427\u00f75=\" },\n  @{ Old = \"200\u00f72=\"; New = \"987\u00f77=\" },\n  @{ Old = \"531\u00f75=\"; New = \"519\u00f75=\" },\n  @{ Old = \"152\u00f77=\"; New = \"223\u00f73=\" },\n  @{ Old = \"670\u00f78=\"; New = \"755\u00f73=\" },\n  @{ Old = \"545\u00f74=\"; New = \"439\u00f75=\" },\n  @{ Old = \"723\u00f72=\"; New = \"688\u00f73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
